# datasync #144: "Test also resamplings."
#   + Add 1 more sheet to be synced (new Sheet4, a copy of the former
#     empty Sheet3 placeholder).
#   + Sheet3 is repurposed to hold a 2x-upsampled ("resampled") copy of
#     Sheet1's x/y1/y2 series, labelled x / y1 / Y_RESAMPLE, and becomes
#     the active sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- grow the tab strip: append a brand-new (empty) Sheet4 after Sheet3 ---
$newSheet = $wb.Worksheets.Add($null, $ws3)

# --- populate Sheet3 with the resampled series -----------------------
$ws3.Range("A1").Value = "x"
$ws3.Range("B1").Value = "y1"
$ws3.Range("C1").Value = "Y_RESAMPLE"

for ($i = 0; $i -lt 59; $i++) {
    $r = $i + 2
    $a = $i / 2.0
    $ws3.Cells.Item($r, 1).Value = $a
    $ws3.Cells.Item($r, 2).Value = $a
    $ws3.Cells.Item($r, 3).Value = 2 * $a
}

# --- selections -------------------------------------------------------
# Sheet1 keeps its data but loses tab focus; its whole used column A
# ends up selected.
$ws1.Range("A1:A31").Select()

# Sheet3 becomes the active sheet/tab with B1 selected.
$ws3.Range("B1").Select()
